# Weekly data refresh: insert three new "Pepino ensalada" price records
# (most recent week) above the existing history, pushing the prior rows
# (920:937) down to (923:940).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 920, shifting the existing data (and its
# formatting, e.g. the date style on column D) down to 923:940.
$ws.Rows("920:922").Insert()

# Row 920 - new record
$ws.Range("A920").Value = 10
$ws.Range("B920").Value = 'Vega Modelo de Temuco'
$ws.Range("C920").Value = 'La Araucanía'
$ws.Range("D920").Value = 45239
$ws.Range("E920").Value = 9
$ws.Range("F920").Value = 100112043
$ws.Range("G920").Value = 'Pepino ensalada'
$ws.Range("H920").Value = 'Sin especificar'
$ws.Range("I920").Value = 'Primera'
$ws.Range("J920").Value = 550
$ws.Range("K920").Value = 16000
$ws.Range("L920").Value = 17000
$ws.Range("M920").Value = 16364
$ws.Range("N920").Value = '$/caja 50 unidades'
$ws.Range("O920").Value = 'Región de Arica y Parinacota'
$ws.Range("P920").Value = 327
$ws.Range("Q920").Value = 50
$ws.Range("R920").Value = 'Hortaliza'

# Row 921 - new record
$ws.Range("A921").Value = 10
$ws.Range("B921").Value = 'Vega Modelo de Temuco'
$ws.Range("C921").Value = 'La Araucanía'
$ws.Range("D921").Value = 45239
$ws.Range("E921").Value = 9
$ws.Range("F921").Value = 100112043
$ws.Range("G921").Value = 'Pepino ensalada'
$ws.Range("H921").Value = 'Sin especificar'
$ws.Range("I921").Value = 'Primera'
$ws.Range("J921").Value = 100
$ws.Range("K921").Value = 19000
$ws.Range("L921").Value = 19000
$ws.Range("M921").Value = 19000
$ws.Range("N921").Value = '$/caja 50 unidades'
$ws.Range("O921").Value = 'Región del Maule'
$ws.Range("P921").Value = 380
$ws.Range("Q921").Value = 50
$ws.Range("R921").Value = 'Hortaliza'

# Row 922 - new record
$ws.Range("A922").Value = 10
$ws.Range("B922").Value = 'Vega Modelo de Temuco'
$ws.Range("C922").Value = 'La Araucanía'
$ws.Range("D922").Value = 45239
$ws.Range("E922").Value = 9
$ws.Range("F922").Value = 100112043
$ws.Range("G922").Value = 'Pepino ensalada'
$ws.Range("H922").Value = 'Sin especificar'
$ws.Range("I922").Value = 'Segunda'
$ws.Range("J922").Value = 50
$ws.Range("K922").Value = 12000
$ws.Range("L922").Value = 13000
$ws.Range("M922").Value = 12600
$ws.Range("N922").Value = '$/caja 60 unidades'
$ws.Range("O922").Value = 'Región del Maule'
$ws.Range("P922").Value = 210
$ws.Range("Q922").Value = 60
$ws.Range("R922").Value = 'Hortaliza'
